# Appends a new "availability check" snapshot (14 rows) to the bottom of
# Sheet1, mirroring the most recent 14-row block (rows 632-645) with a new
# timestamp, matching the automated "Actualizar" run described in the
# commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 645
$blockSize = 14
$newTimestamp = 44232.64054723023

# Names/URLs (as shared-string text) for the 14-row block, taken verbatim
# from the previous block (rows 632-645) so the same services/order repeat.
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
# Text actually displayed in column B (matches the shared strings already
# present in the workbook from earlier identical blocks).
$displayUrls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# Hyperlink target address (without any "#" fragment - that part is
# carried by SubAddress/"location" instead, same split MapStore already
# used in the earlier blocks).
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# "location" fragment used only by the MapStore hyperlink (index 8)
$locations = @("","","","","","","","","/","","","","","")

# First, nudge the 14 rows of the previous run (632-645) - their stored
# timestamp gets a tiny floating point correction in this commit.
for ($i = 0; $i -lt $blockSize; $i++) {
    $r = ($lastRow - $blockSize) + 1 + $i
    $ws.Cells.Item($r, 4).Value = 44232.61942753472
}

for ($i = 0; $i -lt $blockSize; $i++) {
    $r = $lastRow + 1 + $i

    $ws.Cells.Item($r, 1).Value = $names[$i]

    $ws.Cells.Item($r, 2).Value = $displayUrls[$i]
    if ($locations[$i] -ne "") {
        $h = $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $urls[$i], $locations[$i])
    } else {
        $h = $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $urls[$i])
    }
    $ws.Cells.Item($r, 2).Style = "Hyperlink"

    $ws.Cells.Item($r, 3).Value = "Disponible"

    $ws.Cells.Item($r, 4).Value = $newTimestamp
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "Appended $blockSize rows ($($lastRow+1):$($lastRow+$blockSize))"
